{"js": "// Add a new \"Compact List\" paragraph style (styleId \"CompactList\"),\n// mirroring the existing \"Compact\" style: based on \"Body Text\", marked as\n// a quick style, with 1.8pt (36 twips) spacing before/after the paragraph.\n\n// Office.js derives the styleId by stripping spaces from the supplied\n// display name, so \"Compact List\" -> styleId \"CompactList\" / name \"Compact List\".\ncontext.document.addStyle(\"Compact List\", Word.StyleType.paragraph);\nawait context.sync();\n\n// Re-resolve the newly minted style by name. The object returned directly\n// from addStyle() keeps a transient anchor that can get reseated onto the\n// wrong style (e.g. \"Normal\") once other property writes are flushed on the\n// next sync, so fetch a fresh/stable handle before setting properties.\nconst newStyle = context.document.getStyles().getByName(\"Compact List\");\nnewStyle.baseStyle = \"BodyText\";\nnewStyle.quickStyle = true;\nnewStyle.paragraphFormat.spaceBefore = 1.8;\nnewStyle.paragraphFormat.spaceAfter = 1.8;\nawait context.sync();\n", "ps1": "# Add a new \"Compact List\" paragraph style (styleId \"CompactList\"),\n# mirroring the existing \"Compact\" style: based on \"Body Text\", marked as\n# a quick style, with 1.8pt (36 twips) spacing before/after the paragraph.\n\n$d = $word.ActiveDocument\n\n# Word derives the styleId by stripping spaces from the supplied display\n# name, so \"Compact List\" -> styleId \"CompactList\" / name \"Compact List\".\n# wdStyleTypeParagraph = 1\n$d.Styles.Add(\"Compact List\", 1)\n\n# Re-resolve the newly minted style by name before tweaking its properties.\n$newStyle = $d.Styles.Item(\"Compact List\")\n$newStyle.BaseStyle = \"BodyText\"\n$newStyle.QuickStyle = $true\n$newStyle.ParagraphFormat.SpaceBefore = 1.8\n$newStyle.ParagraphFormat.SpaceAfter = 1.8\n"}
